$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "2023" column (T) mirroring the existing "2022" column (S) layout/style.
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 4.8187602774004432
$ws.Range("T6").Value = 11.788953009068425
$ws.Range("T7").Value = 5.2855407047387608
$ws.Range("T8").Value = 11.35112240576027
$ws.Range("T9").Value = 16.577540106951872
$ws.Range("T10").Value = 14.651002073255009
$ws.Range("T11").Value = 5.034965034965035
$ws.Range("T12").Value = 3.1837160751565765
$ws.Range("T13").Value = 2.2263731825525039
$ws.Range("T14").Value = 5.1321450522433931

# Mirror the column S styling for each row into column T.
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("S6:S14").Copy()
$ws.Range("T6:T14").PasteSpecial(-4122)

# Clear the stale cell selection left over from editing (matches the saved view).
$ws.Range("A1").Select()
